$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E3").Value = '[''Normal'']'
$ws.Range("D8").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E8").Value = '[''Normal'', ''HardwareFault'']'
$ws.Range("D12").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E12").Value = '[''Normal'']'
$ws.Range("D15").Value = '[0, 0, 0, 1, 0, 0, 0]'
$ws.Range("E15").Value = '[''ParamViolation'']'
$ws.Range("D16").Value = '[1, 0, 0, 0, 1, 0, 0]'
$ws.Range("E16").Value = '[''Normal'', ''RegulationViolation'']'
$ws.Range("D24").Value = '[0, 0, 0, 0, 0, 0, 0]'
$ws.Range("E24").Value = '[]'
$ws.Range("D27").Value = '[0, 0, 0, 0, 0, 0, 1]'
$ws.Range("E27").Value = '[''SoftwareFault'']'
$ws.Range("D28").Value = '[0, 0, 0, 0, 0, 0, 1]'
$ws.Range("E28").Value = '[''SoftwareFault'']'
$ws.Range("D36").Value = '[1, 1, 1, 0, 0, 0, 0]'
$ws.Range("E36").Value = '[''Normal'', ''SurroundingEnvironment'', ''HardwareFault'']'
$ws.Range("D38").Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Range("E38").Value = '[''Normal'', ''SoftwareFault'']'
$ws.Range("D39").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E39").Value = '[''Normal'']'
$ws.Range("D54").Value = '[0, 0, 0, 0, 0, 1, 0]'
$ws.Range("E54").Value = '[''CommunicationIssue'']'
$ws.Range("D56").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E56").Value = '[''Normal'']'
$ws.Range("D61").Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Range("E61").Value = '[''Normal'']'
$ws.Range("D68").Value = '[1, 0, 0, 1, 0, 0, 0]'
$ws.Range("E68").Value = '[''Normal'', ''ParamViolation'']'
$ws.Range("D69").Value = '[1, 1, 0, 0, 0, 0, 0]'
$ws.Range("E69").Value = '[''Normal'', ''SurroundingEnvironment'']'
$ws.Range("D80").Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Range("E80").Value = '[''Normal'', ''HardwareFault'']'
$ws.Range("D92").Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Range("E92").Value = '[''Normal'', ''SoftwareFault'']'
$ws.Range("D93").Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Range("E93").Value = '[''Normal'', ''SoftwareFault'']'
$ws.Range("D109").Value = '[1, 1, 0, 0, 0, 0, 0]'
$ws.Range("E109").Value = '[''Normal'', ''SurroundingEnvironment'']'
